$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New strings are added in this order to reproduce the shared-strings table order:
# 8  Tracking data visualisation   -> A8
# 9  rmd file name                 -> E1
# 10 03-TrackingData-Visualisation -> E8
# 11 track-vis                     -> D8
# 12 introduction                  -> B3 (replaces old "tbc")
# 13 swe1                          -> D4
# 14 swe2                          -> D5
# 15 swe3                          -> D6
# 16 single colony                 -> C5
# 17 multiple colonies             -> C7
# 18 background raster single      -> C4
# 19 background raster multi       -> C6
# 20 swe4                          -> D7

$ws.Range("A8").Value = "Tracking data visualisation"
$ws.Range("E1").Value = "rmd file name"
$ws.Range("E8").Value = "03-TrackingData-Visualisation"
$ws.Range("D8").Value = "track-vis"
$ws.Range("B3").Value = "introduction"
$ws.Range("D4").Value = "swe1"
$ws.Range("D5").Value = "swe2"
$ws.Range("D6").Value = "swe3"
$ws.Range("C5").Value = "single colony"
$ws.Range("C7").Value = "multiple colonies"
$ws.Range("C4").Value = "background raster single"
$ws.Range("C6").Value = "background raster multi"
$ws.Range("D7").Value = "swe4"

# Remove the old C3 value ("tbc") which is no longer present in the new layout
$ws.Range("C3").ClearContents()

# Remaining row 8 cells that reuse existing shared strings
$ws.Range("B8").Value = "tbc"
$ws.Range("C8").Value = "tbc"

# Column widths (values chosen to round-trip as closely as possible to the
# target stored widths of 22.81640625 / 12.90625 characters)
$ws.Columns.Item(3).ColumnWidth = 22
$ws.Columns.Item(5).ColumnWidth = 12

# Update the active selection to match the edited file
$ws.Range("D4").Select() | Out-Null
